$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are numeric-looking strings that must stay as text.
# Prefix with a literal leading apostrophe (Excel quote-prefix) to force text entry,
# then reset the cell style back to Normal so no extra formatting/style is introduced.

$ws.Range("D2").Value = "'37.833.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.05%  "

$ws.Range("D3").Value = "'2.079.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.15%  "

$ws.Range("D5").Value = "'233.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("E6").Value = "  +0.42%  "

$ws.Range("D7").Value = "'59.43"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.91%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("E10").Value = "  +1.41%  "

$ws.Range("E11").Value = "  +0.85%  "

$ws.Range("D12").Value = "'14.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.89%  "

$ws.Range("D13").Value = "'21.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.18%  "

$ws.Range("D14").Value = "'0.776"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.31%  "

$ws.Range("D15").Value = "'5.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.08%  "

$ws.Range("D16").Value = "'2.078.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.13%  "

$ws.Range("D17").Value = "'37.727.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("D18").Value = "'6.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("D19").Value = "'71.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.75%  "

$ws.Range("D20").Value = "'0.0₃0852"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.28%  "

$ws.Range("D21").Value = "'228.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").Value = "'2.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.25%  "

$ws.Range("D24").Value = "'2.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.00%  "

$ws.Range("D25").Value = "'170.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.60%  "

$ws.Range("D26").Value = "'9.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.28%  "

$ws.Range("E27").Value = "  -2.65%  "

$ws.Range("E28").Value = "  -0.21%  "

$ws.Range("D29").Value = "'19.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  +1.75%  "

$ws.Range("D31").Value = "'4.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.24%  "

$ws.Range("D32").Value = "'4.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.74%  "

$ws.Range("D33").Value = "'0.0634"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.51%  "

$ws.Range("D34").Value = "'2.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.93%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("E36").Value = "  -0.75%  "

$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("D38").Value = "'5.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.36%  "

$ws.Range("D39").Value = "'0.0977"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.49%  "

$ws.Range("D40").Value = "'99.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.88%  "

$ws.Range("D41").Value = "'0.0217"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.57%  "

$ws.Range("E42").Value = "  -1.78%  "

$ws.Range("D43").Value = "'16.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.47%  "

$ws.Range("D44").Value = "'1.445.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.80%  "

$ws.Range("D45").Value = "'1.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.90%  "

$ws.Range("D46").Value = "'4.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.81%  "

$ws.Range("E47").Value = "  +0.33%  "

$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("E49").Value = "  +0.27%  "

$ws.Range("D50").Value = "'2.269.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.31%  "

$ws.Range("D51").Value = "'46.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.70%  "
